$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='243.52' },
    @{ Row=3; D='25.01' },
    @{ Row=4; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='3.498'; E='3LEOLEO' },
    @{ Row=5; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='5.157'; E='4HuobiTokenHT' },
    @{ Row=6; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.05728'; E='5CronosCRO' },
    @{ Row=7; B='KuCoinToken'; C='https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'; D='6.482'; E='6KuCoinTokenKCS' },
    @{ Row=8; B='GateToken'; C='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D='3.078'; E='7GateTokenGT' },
    @{ Row=9; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='0.8101'; E='8MXTokenMX' },
    @{ Row=10; B='FTXToken'; C='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; D='0.8408'; E='9FTXTokenFTT' },
    @{ Row=11; B='WazirX'; C='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D='0.1337'; E='10WazirXWRX' },
    @{ Row=12; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.06954'; E='11MandalaExchangeTokenMDX' },
    @{ Row=13; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.02823'; E='12BitrueCoinBTR' },
    @{ Row=14; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.09362'; E='13BitMartTokenBMX' },
    @{ Row=15; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001509'; E='14BitForexTokenBF' },
    @{ Row=16; B='One'; C='https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; D='0.0005986'; E='15OneONE' },
    @{ Row=17; B='TigerCash'; C='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D='0.006109'; E='16TigerCashTCH' },
    @{ Row=18; D='2.120' },
    @{ Row=19; D='0.3195' },
    @{ Row=20; D='0.03149' },
    @{ Row=21; D='0.1300' },
    @{ Row=22; D='3.741' },
    @{ Row=23; D='0.04669' },
    @{ Row=24; D='0.1328' },
    @{ Row=25; D='0.001239' },
    @{ Row=26; D='0.004262' },
    @{ Row=27; D='0.00009697' },
    @{ Row=28; D='0.0001500'; E='27UpBotsUBXTWorstin24h' },
    @{ Row=40; D='0.03616' },
    @{ Row=41; D='0.006288' },
    @{ Row=42; D='0.1048' },
    @{ Row=43; D='0.002999' },
    @{ Row=44; D='0.007283' },
    @{ Row=45; D='0.00005267' },
    @{ Row=47; D='0.1998' },
    @{ Row=49; D='0.00002099' },
    @{ Row=50; D='0.0001999' }
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Range("B$r").Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Range("C$r").Value = $item.C }
    if ($item.ContainsKey("D")) {
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.ClearFormats()
    }
    if ($item.ContainsKey("E")) { $ws.Range("E$r").Value = $item.E }
}
